# Adds a new paragraph "Und noch ein test" right after the first
# paragraph ("Dies ist ein Test"), keeping the pre-existing trailing
# empty paragraph intact. "test" is wrapped in spellStart/spellEnd
# proofErr markers, the way Word's background spell checker flags it
# while typing.

$d = $word.ActiveDocument

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# Locate the first paragraph ("Dies ist ein Test") and the paragraph
# that currently follows it (the trailing empty paragraph already
# present in the document).
$firstPara = $d.Paragraphs(1)
$trailingPara = $d.Paragraphs(2)

# Capture the trailing paragraph's own OOXML so we can preserve its
# identity (paraId/rsid/etc.) once we rewrite this span.
$trailingXml = $trailingPara.Range.WordOpenXML
if ($trailingXml -match '(<w:p\b[^>]*/>)') {
    $trailingParaXml = $matches[1]
} elseif ($trailingXml -match '(<w:p\b[^>]*>.*?</w:p>)') {
    $trailingParaXml = $matches[1]
} else {
    $trailingParaXml = "<w:p/>"
}

# New paragraph's OOXML: "Und noch ein " + proofErr-wrapped "test".
$newParaXml = '<w:p ' + $wNs + '>' +
    '<w:r><w:t xml:space="preserve">Und noch ein </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>test</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '</w:p>'

# Replace the span from the end of the first paragraph through the end
# of the trailing paragraph with [new paragraph][trailing paragraph].
# InsertXML replaces exactly the range it's called on, so by scoping it
# to this paragraph-aligned span we leave "Dies ist ein Test" and the
# document's sectPr untouched.
$span = $d.Range($firstPara.Range.End, $trailingPara.Range.End)
$span.InsertXML($newParaXml + $trailingParaXml)
